# Update task schedule values (used in testing)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where y_corrSteps (E), y_nrSteps (G) and alienID (H) values changed
$rows = @(4, 8, 16, 18, 23, 27)
$eVals = @{ 4 = 5; 8 = 5; 16 = 6; 18 = 5; 23 = 4; 27 = 6 }
$gVals = @{ 4 = -3; 8 = -3; 16 = -3; 18 = -3; 23 = -3; 27 = -3 }
$hVals = @{ 4 = 13; 8 = 13; 16 = 13; 18 = 13; 23 = 13; 27 = 13 }

foreach ($r in $rows) {
    $ws.Range("E$r").Value = $eVals[$r]
    $ws.Range("G$r").Value = $gVals[$r]
    $ws.Range("H$r").Value = $hVals[$r]
}
